$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column widths (target width = ColumnWidth input + 0.8333333333333334 offset)
$ws.Columns.Item(2).ColumnWidth = 7.166666666666667
$ws.Columns.Item(3).ColumnWidth = 7.166666666666667
$ws.Columns.Item(6).ColumnWidth = 7.166666666666667
$ws.Columns.Item(7).ColumnWidth = 7.166666666666667
$ws.Columns.Item(9).ColumnWidth = 7.166666666666667
$ws.Columns.Item(10).ColumnWidth = 7.166666666666667
$ws.Columns.Item(11).ColumnWidth = 7.166666666666667
$ws.Columns.Item(12).ColumnWidth = 7.166666666666667
$ws.Columns.Item(13).ColumnWidth = 7.166666666666667
$ws.Columns.Item(15).ColumnWidth = 7.166666666666667
$ws.Columns.Item(16).ColumnWidth = 7.166666666666667
$ws.Columns.Item(17).ColumnWidth = 7.166666666666667
$ws.Columns.Item(20).ColumnWidth = 8.166666666666666
$ws.Columns.Item(22).ColumnWidth = 7.166666666666667
$ws.Columns.Item(23).ColumnWidth = 7.166666666666667
$ws.Columns.Item(24).ColumnWidth = 7.166666666666667
$ws.Columns.Item(26).ColumnWidth = 7.166666666666667
$ws.Columns.Item(27).ColumnWidth = 7.166666666666667
$ws.Columns.Item(28).ColumnWidth = 7.166666666666667
$ws.Columns.Item(29).ColumnWidth = 7.166666666666667
$ws.Columns.Item(30).ColumnWidth = 7.166666666666667
$ws.Columns.Item(34).ColumnWidth = 7.166666666666667

# Update cell values for rows 2-5 (new data)
# Row 2
$ws.Cells.Item(2, 1).Value2 = 45080.50694444445
$ws.Cells.Item(2, 2).Value2 = 21.139
$ws.Cells.Item(2, 3).Value2 = 14.405
$ws.Cells.Item(2, 4).Value2 = 4.093
$ws.Cells.Item(2, 5).Value2 = 44.473
$ws.Cells.Item(2, 6).Value2 = 36.592
$ws.Cells.Item(2, 7).Value2 = 16.635
$ws.Cells.Item(2, 8).Value2 = 54.157
$ws.Cells.Item(2, 9).Value2 = 25.596
$ws.Cells.Item(2, 10).Value2 = 10.793
$ws.Cells.Item(2, 11).Value2 = 16.568
$ws.Cells.Item(2, 12).Value2 = 17.656
$ws.Cells.Item(2, 13).Value2 = 18.451
$ws.Cells.Item(2, 14).Value2 = 5.311
$ws.Cells.Item(2, 15).Value2 = 16.542
$ws.Cells.Item(2, 16).Value2 = 23.194
$ws.Cells.Item(2, 17).Value2 = 13.966
$ws.Cells.Item(2, 18).Value2 = 3.726
$ws.Cells.Item(2, 19).Value2 = 2.455
$ws.Cells.Item(2, 20).Value2 = 244.445
$ws.Cells.Item(2, 21).Value2 = 46.035
$ws.Cells.Item(2, 22).Value2 = 15.269
$ws.Cells.Item(2, 23).Value2 = 30.403
$ws.Cells.Item(2, 24).Value2 = 15.788
$ws.Cells.Item(2, 25).Value2 = 2.61
$ws.Cells.Item(2, 26).Value2 = 26.878
$ws.Cells.Item(2, 27).Value2 = 13.487
$ws.Cells.Item(2, 28).Value2 = 12.177
$ws.Cells.Item(2, 29).Value2 = 14.251
$ws.Cells.Item(2, 30).Value2 = 18.232
$ws.Cells.Item(2, 31).Value2 = 3.64
$ws.Cells.Item(2, 32).Value2 = 47.882
$ws.Cells.Item(2, 33).Value2 = 8.427
$ws.Cells.Item(2, 34).Value2 = 19.09
# Row 3
$ws.Cells.Item(3, 1).Value2 = 45080.51388888889
$ws.Cells.Item(3, 2).Value2 = 12.011
$ws.Cells.Item(3, 3).Value2 = 8.323
$ws.Cells.Item(3, 4).Value2 = 1.66
$ws.Cells.Item(3, 5).Value2 = 25.537
$ws.Cells.Item(3, 6).Value2 = 21.012
$ws.Cells.Item(3, 7).Value2 = 9.452
$ws.Cells.Item(3, 8).Value2 = 38.489
$ws.Cells.Item(3, 9).Value2 = 14.543
$ws.Cells.Item(3, 10).Value2 = 6.213
$ws.Cells.Item(3, 11).Value2 = 9.34
$ws.Cells.Item(3, 12).Value2 = 10.304
$ws.Cells.Item(3, 13).Value2 = 10.705
$ws.Cells.Item(3, 14).Value2 = 3.021
$ws.Cells.Item(3, 15).Value2 = 9.398999999999999
$ws.Cells.Item(3, 16).Value2 = 13.209
$ws.Cells.Item(3, 17).Value2 = 8.148999999999999
$ws.Cells.Item(3, 18).Value2 = 1.571
$ws.Cells.Item(3, 19).Value2 = 0.945
$ws.Cells.Item(3, 20).Value2 = 135.752
$ws.Cells.Item(3, 21).Value2 = 26.38
$ws.Cells.Item(3, 22).Value2 = 8.676
$ws.Cells.Item(3, 23).Value2 = 17.377
$ws.Cells.Item(3, 24).Value2 = 9.237
$ws.Cells.Item(3, 25).Value2 = 1.558
$ws.Cells.Item(3, 26).Value2 = 18.183
$ws.Cells.Item(3, 27).Value2 = 7.663
$ws.Cells.Item(3, 28).Value2 = 6.995
$ws.Cells.Item(3, 29).Value2 = 8.182
$ws.Cells.Item(3, 30).Value2 = 10.692
$ws.Cells.Item(3, 31).Value2 = 1.294
$ws.Cells.Item(3, 32).Value2 = 34.997
$ws.Cells.Item(3, 33).Value2 = 4.75
$ws.Cells.Item(3, 34).Value2 = 10.847
# Row 4
$ws.Cells.Item(4, 1).Value2 = 45080.52083333334
$ws.Cells.Item(4, 2).Value2 = 0.379
$ws.Cells.Item(4, 3).Value2 = 0.173
$ws.Cells.Item(4, 4).Value2 = 0.786
$ws.Cells.Item(4, 5).Value2 = 0.679
$ws.Cells.Item(4, 6).Value2 = 0.192
$ws.Cells.Item(4, 7).Value2 = 0
$ws.Cells.Item(4, 8).Value2 = 8.196999999999999
$ws.Cells.Item(4, 9).Value2 = 0.582
$ws.Cells.Item(4, 10).Value2 = 0.248
$ws.Cells.Item(4, 11).Value2 = 0.066
$ws.Cells.Item(4, 12).Value2 = 0.303
$ws.Cells.Item(4, 13).Value2 = 0
$ws.Cells.Item(4, 14).Value2 = 0
$ws.Cells.Item(4, 15).Value2 = 0.376
$ws.Cells.Item(4, 16).Value2 = 0.614
$ws.Cells.Item(4, 17).Value2 = 0.613
$ws.Cells.Item(4, 18).Value2 = 0.9399999999999999
$ws.Cells.Item(4, 19).Value2 = 0.309
$ws.Cells.Item(4, 20).Value2 = 0
$ws.Cells.Item(4, 21).Value2 = 1.632
$ws.Cells.Item(4, 22).Value2 = 0.347
$ws.Cells.Item(4, 23).Value2 = 1.022
$ws.Cells.Item(4, 24).Value2 = 0.5580000000000001
$ws.Cells.Item(4, 25).Value2 = 0.309
$ws.Cells.Item(4, 26).Value2 = 3.283
$ws.Cells.Item(4, 27).Value2 = 0.307
$ws.Cells.Item(4, 28).Value2 = 0.407
$ws.Cells.Item(4, 29).Value2 = 0.426
$ws.Cells.Item(4, 30).Value2 = 0.281
$ws.Cells.Item(4, 31).Value2 = 0.783
$ws.Cells.Item(4, 32).Value2 = 7.976
$ws.Cells.Item(4, 33).Value2 = 0.046
$ws.Cells.Item(4, 34).Value2 = 0.46
# Row 5
$ws.Cells.Item(5, 1).Value2 = 45080.52777777778
$ws.Cells.Item(5, 2).Value2 = 8.56
$ws.Cells.Item(5, 3).Value2 = 6.35
$ws.Cells.Item(5, 4).Value2 = 0.84
$ws.Cells.Item(5, 5).Value2 = 18.55
$ws.Cells.Item(5, 6).Value2 = 15.04
$ws.Cells.Item(5, 7).Value2 = 6.97
$ws.Cells.Item(5, 8).Value2 = 23.23
$ws.Cells.Item(5, 9).Value2 = 10.47
$ws.Cells.Item(5, 10).Value2 = 4.5
$ws.Cells.Item(5, 11).Value2 = 6.76
$ws.Cells.Item(5, 12).Value2 = 7.47
$ws.Cells.Item(5, 13).Value2 = 7.69
$ws.Cells.Item(5, 14).Value2 = 2.07
$ws.Cells.Item(5, 15).Value2 = 6.77
$ws.Cells.Item(5, 16).Value2 = 9.470000000000001
$ws.Cells.Item(5, 17).Value2 = 5.86
$ws.Cells.Item(5, 18).Value2 = 0.76
$ws.Cells.Item(5, 19).Value2 = 0.46
$ws.Cells.Item(5, 20).Value2 = 96.06999999999999
$ws.Cells.Item(5, 21).Value2 = 18.7
$ws.Cells.Item(5, 22).Value2 = 6.25
$ws.Cells.Item(5, 23).Value2 = 12.35
$ws.Cells.Item(5, 24).Value2 = 6.69
$ws.Cells.Item(5, 25).Value2 = 1.1
$ws.Cells.Item(5, 26).Value2 = 11.14
$ws.Cells.Item(5, 27).Value2 = 5.52
$ws.Cells.Item(5, 28).Value2 = 4.94
$ws.Cells.Item(5, 29).Value2 = 5.8
$ws.Cells.Item(5, 30).Value2 = 7.84
$ws.Cells.Item(5, 31).Value2 = 0.5600000000000001
$ws.Cells.Item(5, 32).Value2 = 20.63
$ws.Cells.Item(5, 33).Value2 = 3.44
$ws.Cells.Item(5, 34).Value2 = 7.81

# Delete row 6 (no longer present in target data)
$ws.Rows.Item(6).Delete()
